# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows (Cereza - Brooks/Lapins/Rainier, Paine/Curicó,
# fecha serial 44914) right before the existing row 32, pushing the previous
# data down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above row 32; Excel copies formatting (incl. the
# date number format on column D) from the row that was at 32 down into the
# newly-created rows, same as a manual "Insert Copied Cells" / row insert.
$ws.Rows("32:34").Insert()

# --- New row 32: Cereza / Brooks / Primera ---
$ws.Range("A32").Value = 4
$ws.Range("B32").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C32").Value = "Los Lagos"
$ws.Range("D32").Value = 44914
$ws.Range("E32").Value = 10
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100103
$ws.Range("H32").Value = "Frutos de hueso (carozo)"
$ws.Range("I32").Value = 100103001
$ws.Range("J32").Value = "Cereza"
$ws.Range("K32").Value = "Brooks"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 800
$ws.Range("N32").Value = 6500
$ws.Range("O32").Value = 7000
$ws.Range("P32").Value = 6750
$ws.Range("Q32").Value = "$/bandeja 10 kilos"
$ws.Range("R32").Value = "Provincia de Curicó"
$ws.Range("S32").Value = 675
$ws.Range("T32").Value = 10

# --- New row 33: Cereza / Lapins / Primera ---
$ws.Range("A33").Value = 4
$ws.Range("B33").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C33").Value = "Los Lagos"
$ws.Range("D33").Value = 44914
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100103
$ws.Range("H33").Value = "Frutos de hueso (carozo)"
$ws.Range("I33").Value = 100103001
$ws.Range("J33").Value = "Cereza"
$ws.Range("K33").Value = "Lapins"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 800
$ws.Range("N33").Value = 6500
$ws.Range("O33").Value = 7000
$ws.Range("P33").Value = 6750
$ws.Range("Q33").Value = "$/bandeja 10 kilos"
$ws.Range("R33").Value = "Provincia de Curicó"
$ws.Range("S33").Value = 675
$ws.Range("T33").Value = 10

# --- New row 34: Cereza / Rainier / Primera ---
$ws.Range("A34").Value = 4
$ws.Range("B34").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C34").Value = "Los Lagos"
$ws.Range("D34").Value = 44914
$ws.Range("E34").Value = 10
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100103
$ws.Range("H34").Value = "Frutos de hueso (carozo)"
$ws.Range("I34").Value = 100103001
$ws.Range("J34").Value = "Cereza"
$ws.Range("K34").Value = "Rainier"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 400
$ws.Range("N34").Value = 9000
$ws.Range("O34").Value = 10000
$ws.Range("P34").Value = 9500
$ws.Range("Q34").Value = "$/bandeja 10 kilos"
$ws.Range("R34").Value = "Paine"
$ws.Range("S34").Value = 950
$ws.Range("T34").Value = 10
